# Rename the IEA data-file references (now plain .csv names instead of the
# old .xlsx naming) and add the new Australia ("AU.csv") entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("country_file_name")

# Italy row (row 2): country data iea -> "IT.csv"
$ws.Range("C2").Value = "IT.csv"

# Egypt row (row 3): country data iea -> "EG.csv"
$ws.Range("C3").Value = "EG.csv"

# Australia row (row 6): new "country data iea" entry -> "AU.csv"
$ws.Range("C6").Value = "AU.csv"

# Match the author's final cursor/selection position recorded in the file.
[void]$ws.Range("G11").Select()
